# Adds columns I ("I0") and J ("IF") with header + 39 rows of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (H1) onto the two new header cells before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(8, 9),
    @(11, 12),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(5, 6),
    @(6, 6),
    @(7, 8),
    @(4, 4),
    @(6, 6),
    @(5, 5),
    @(7, 7),
    @(8, 9),
    @(6, 6),
    @(4, 4),
    @(3, 3),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
